# Update cryptos list values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.857.98"
$ws.Range("E2").Value = "  +3.14%  "

$ws.Range("D3").Value = "3.453.37"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.19"
$ws.Range("E5").Value = "  +2.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.22"
$ws.Range("E6").Value = "  +5.07%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.478"
$ws.Range("E8").Value = "  +1.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.64"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  +2.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.392"
$ws.Range("E11").Value = "  +1.32%  "

$ws.Range("D12").Value = "4.013.13"
$ws.Range("E12").Value = "  +1.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.20"
$ws.Range("E13").Value = "  +5.38%  "

$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").Value = "3.423.53"
$ws.Range("E15").Value = "  +0.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  +2.08%  "

$ws.Range("D17").Value = "60.876.99"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.24"
$ws.Range("E18").Value = "  +2.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.24"
$ws.Range("E19").Value = "  +5.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.30"
$ws.Range("E20").Value = "  +4.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "396.84"
$ws.Range("E21").Value = "  +4.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.68"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.561"
$ws.Range("E23").Value = "  +2.31%  "

$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000117"
$ws.Range("E25").Value = "  +3.14%  "

$ws.Range("D26").Value = "3.559.46"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.189"
$ws.Range("E27").Value = "  -1.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.64"
$ws.Range("E28").Value = "  +6.17%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.14"
$ws.Range("E30").Value = "  +2.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.15"
$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  +4.17%  "

$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.86"
$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  +8.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.05"
$ws.Range("E36").Value = "  +2.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.59"
$ws.Range("E37").Value = "  +9.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.80"
$ws.Range("E38").Value = "  +1.26%  "

$ws.Range("D39").Value = "3.459.03"
$ws.Range("E39").Value = "  +1.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "29.38"
$ws.Range("E40").Value = "  +13.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0765"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.793"
$ws.Range("E42").Value = "  +1.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.48"
$ws.Range("E43").Value = "  +3.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.72"
$ws.Range("E44").Value = "  +5.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.19"
$ws.Range("E45").Value = "  +7.00%  "

$ws.Range("D46").Value = "2.522.82"
$ws.Range("E46").Value = "  +3.83%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.22"
$ws.Range("E47").Value = "  +2.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.71"
$ws.Range("E48").Value = "  +1.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.995"
$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0267"
$ws.Range("E50").Value = "  +2.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.14"
$ws.Range("E51").Value = "  +2.83%  "
